$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.345.07'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '2.389.08'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '520.85'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.97'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.50'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.86%  '
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.343'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.46'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.14%  '
$ws.Range('D14').Value = '2.788.46'
$ws.Range('E14').Value = '  +1.26%  '
$ws.Range('D15').Value = '57.340.58'
$ws.Range('E15').Value = '  +0.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000135'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').Value = '2.367.52'
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.59'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '330.08'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.74%  '
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.48'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('E24').Value = '  +16.15%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.165'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.37%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.35'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +11.93%  '
$ws.Range('D28').Value = '0.0₃0751'
$ws.Range('E28').Value = '  +1.38%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.70'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.91%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '166.70'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.01%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.29'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.63'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.24%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  +3.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.996'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.922'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -4.55%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  +5.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '38.76'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '150.07'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.387'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '292.35'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.97%  '
$ws.Range('E43').Value = '  +1.66%  '
$ws.Range('E44').Value = '  +2.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0941'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('E47').Value = '  +0.71%  '
$ws.Range('E48').Value = '  +5.15%  '
$ws.Range('E49').Value = '  +1.53%  '
$ws.Range('E50').Value = '  +4.24%  '
$ws.Range('B51').Value = 'Polygon'
$ws.Range('C51').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.362'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -5.48%  '
